$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = -0.3864313882071876
$ws.Range("J3").Value = 0.2088418595083264
$ws.Range("K3").Value = -0.7052803262679067
$ws.Range("L3").Value = 3.01479743913432
